$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 20 (2025Q2) metrics: total_customers, returning_customers, new_customers, recurrence_rate
$ws.Range("C20").Value = 280
$ws.Range("D20").Value = 229
$ws.Range("E20").Value = 51
$ws.Range("F20").Value = 75.32894736842105
